$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values: shift data up, removing old C2 and old C11 values
$ws.Range("C2").Value = 59.3
$ws.Range("C3").Value = 54.1
$ws.Range("C4").Value = 61.6
$ws.Range("C5").Value = 46.3
$ws.Range("C6").Value = 48.1
$ws.Range("C7").Value = 58.5
$ws.Range("C8").Value = 60.3
$ws.Range("C9").Value = 44.1
$ws.Range("C10").Value = 65.1
$ws.Range("C11").Value = 65.7
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()

# Update the active cell selection to D13
$ws.Range("D13").Select()
